$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5988863706588745
$ws.Range("B1").Value = 0.8039505481719971
$ws.Range("C1").Value = 4.249866008758545
$ws.Range("D1").Value = 2.005975723266602
$ws.Range("E1").Value = 1.047523856163025
